$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.176.29"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.588.44"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "211.48"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.0605"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "19.01"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.811.65"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.593.14"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "63.52"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "26.175.25"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "7.39"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "214.28"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "144.18"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "6.96"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "15.05"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "1.417.86"
$ws.Range("E33").Value = "  +8.92%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "0.583"
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "0.946"
$ws.Range("E42").Value = "  -13.56%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.762"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.12"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "1.723.30"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "86.26"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.49"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0500"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0955"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.14%  "
